# Project pr2 SAVE: update the "Integer max" rule value for R10 (row 8)
# in the "Rules" sheet from 11 to 112.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("D8").Value = 112
